## Actualización automática 2025-10-07 10:30:08
## New asesor "TORRES CADENA JAVIER JOSUE" added to OFICINA-CATAECSA group,
## plus a new sale recorded for PALMA PICO OSCAR FILIDEL (PIEDRA SINTERIZADA).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# New sale amount for PALMA PICO OSCAR FILIDEL (row 36, column L = PIEDRA SINTERIZADA)
$ws1.Range("L36").Value = 179.12

# Insert the new asesor row before the existing row 43 (VACA PANCHI DORYS CAROLINA...)
$ws1.Rows.Item(43).Insert()

$ws1.Range("A43").Value = "OFICINA-CATAECSA"
$ws1.Range("B43").Value = "TORRES CADENA JAVIER JOSUE"
$ws1.Range("C43:R43").Value = 0

# Update the trailing "x de N" summary row (now shifted from 48 to 49)
for ($col = 3; $col -le 18; $col++) {
    $ws1.Cells.Item(49, $col).Value = "0 de 47"
}
$ws1.Cells.Item(49, 12).Value = "1 de 47"

# ---------------------------------------------------------------------
# Sheet "VENTA MENSUAL"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# Same new sale amount mirrored here (column F = PIEDRA SINTERIZADA)
$ws2.Range("F36").Value = 179.12

# Insert the new asesor row before the existing row 43
$ws2.Rows.Item(43).Insert()

$ws2.Range("A43").Value = "OFICINA-CATAECSA"
$ws2.Range("B43").Value = "TORRES CADENA JAVIER JOSUE"
$ws2.Range("C43:G43").Value = 0

# Update the totals row (now shifted from 48 to 49) to reflect the new sale
$ws2.Range("C49").Value = 2652.23
$ws2.Range("D49").Value = 1566.5
$ws2.Range("E49").Value = 13412.01
$ws2.Range("F49").Value = 146.39
$ws2.Range("G49").Value = 0

# ---------------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$ws3.Range("D2").Value = 2818.29
$ws3.Range("E2").Value = -2818.29

$ws3.Range("D4").Value = 2833.46
$ws3.Range("E4").Value = 17166.54
$ws3.Range("F4").Value = 0.141673
